# The edit permutes the data rows 2-8 of the "Artfynd" sheet: each row's
# whole record (Id, coordinates, species fields, and the public-comment
# field) moves to a different row position. Row 9 and the header row are
# untouched. We apply the new values directly to each destination cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (becomes what used to be row 5's location/id/comment) ---
$ws.Range("A2").Value = 111639173
$ws.Range("Q2").Value = 547838.0352795018
$ws.Range("R2").Value = 6926228.915831603
$ws.Range("AC2").Value = "ca 15 plantor"

# --- Row 3 (becomes what used to be row 7's location/id) ---
$ws.Range("A3").Value = 111639175
$ws.Range("Q3").Value = 547828.4099300706
$ws.Range("R3").Value = 6926124.660841302

# --- Row 4 (becomes what used to be row 6's location/id; loses its comment) ---
$ws.Range("A4").Value = 111639172
$ws.Range("Q4").Value = 548221.3480213688
$ws.Range("R4").Value = 6926511.607424877
$ws.Range("AC4").ClearContents()

# --- Row 5 (becomes what used to be row 4's location/id/comment) ---
$ws.Range("A5").Value = 111639167
$ws.Range("Q5").Value = 547814.5103353403
$ws.Range("R5").Value = 6926124.461383951
$ws.Range("AC5").Value = "1 planta"

# --- Row 6 (becomes what used to be row 8's location/id/comment) ---
$ws.Range("A6").Value = 111639170
$ws.Range("Q6").Value = 548231.4260436196
$ws.Range("R6").Value = 6926519.619127685
$ws.Range("AC6").Value = "ca 15 plantor"

# --- Row 7 (becomes what used to be row 2's full species record/comment) ---
$ws.Range("A7").Value = 111639169
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("Q7").Value = 548224.5774945696
$ws.Range("R7").Value = 6926512.579557057
$ws.Range("AC7").Value = "riklig förekomst, mer än 50 plantor"

# --- Row 8 (becomes what used to be row 3's full species record; loses comment) ---
$ws.Range("A8").Value = 111639168
$ws.Range("B8").Value = 89686
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 658
$ws.Range("F8").Value = "Rosenticka"
$ws.Range("G8").Value = "Rhodofomes roseus"
$ws.Range("H8").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q8").Value = 548104.1391889038
$ws.Range("R8").Value = 6926477.987023209
$ws.Range("AC8").ClearContents()
